$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("October")

$ws.Range("B2").Value = 1788
$ws.Range("C2").Value = 1239
$ws.Range("D2").Value = 549
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.44 : 1"

$ws.Range("B3").Value = 541
$ws.Range("C3").Value = 510
$ws.Range("D3").Value = 31
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.06 : 1"

$ws.Range("B4").Value = 1106
$ws.Range("C4").Value = 1392
$ws.Range("D4").Value = -286
$ws.Range("F4").Value = "We lent more than we borrowed"
$ws.Range("G4").Value = "0.79 : 1"

$ws.Range("B5").Value = 82
$ws.Range("C5").Value = 144
$ws.Range("D5").Value = -62
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.57 : 1"

$ws.Range("B6").Value = 1249
$ws.Range("C6").Value = 1516
$ws.Range("D6").Value = -267
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.82 : 1"

$ws.Range("B7").Value = 219
$ws.Range("C7").Value = 200
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = "We borrowerd more than we lent"
$ws.Range("G7").Value = "1.10 : 1"

$ws.Range("B8").Value = 81
$ws.Range("C8").Value = 179
$ws.Range("D8").Value = -98
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.45 : 1"

$ws.Range("B9").Value = 52
$ws.Range("C9").Value = 66
$ws.Range("D9").Value = -14
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.79 : 1"

$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 57
$ws.Range("D10").Value = -56
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.02 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 23
$ws.Range("D12").Value = -3
$ws.Range("F12").Value = "We lent more than we borrowed"
$ws.Range("G12").Value = "0.87 : 1"

$ws.Range("B13").Value = 182
$ws.Range("C13").Value = 125
$ws.Range("D13").Value = 57
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.46 : 1"

$ws.Range("B14").Value = 119
$ws.Range("C14").Value = 316
$ws.Range("D14").Value = -197
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.38 : 1"

$ws.Range("B15").Value = 85
$ws.Range("C15").Value = 123
$ws.Range("D15").Value = -38
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.69 : 1"

$ws.Range("B16").Value = 21
$ws.Range("C16").Value = 144
$ws.Range("D16").Value = -123
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.15 : 1"

$ws.Range("B17").Value = 769
$ws.Range("C17").Value = 480
$ws.Range("D17").Value = 289
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.60 : 1"

$ws.Range("B18").Value = 74
$ws.Range("C18").Value = 118
$ws.Range("D18").Value = -44
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.63 : 1"

$ws.Range("B19").Value = 619
$ws.Range("C19").Value = 360
$ws.Range("D19").Value = 259
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.72 : 1"

$ws.Range("B20").Value = 11
$ws.Range("C20").Value = 72
$ws.Range("D20").Value = -61
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.15 : 1"

$ws.Range("B21").Value = 441
$ws.Range("C21").Value = 456
$ws.Range("D21").Value = -15
$ws.Range("F21").Value = "We lent more than we borrowed"
$ws.Range("G21").Value = "0.97 : 1"

$ws.Range("B22").Value = 48
$ws.Range("C22").Value = 108
$ws.Range("D22").Value = -60
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.44 : 1"

$ws.Range("B23").Value = 776
$ws.Range("C23").Value = 414
$ws.Range("D23").Value = 362
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.87 : 1"

$ws.Range("B24").Value = 1710
$ws.Range("C24").Value = 1361
$ws.Range("D24").Value = 349
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.26 : 1"

$ws.Range("B25").Value = 188
$ws.Range("C25").Value = 354
$ws.Range("D25").Value = -166
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.53 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 234
$ws.Range("C27").Value = 219
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.07 : 1"

$ws.Range("B28").Value = 41
$ws.Range("C28").Value = 59
$ws.Range("D28").Value = -18
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.69 : 1"

$ws.Range("B29").Value = 583
$ws.Range("C29").Value = 491
$ws.Range("D29").Value = 92
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.19 : 1"

$ws.Range("B30").Value = 36
$ws.Range("C30").Value = 41
$ws.Range("D30").Value = -5
$ws.Range("F30").Value = "We lent more than we borrowed"
$ws.Range("G30").Value = "0.88 : 1"

$ws.Range("B31").Value = 116
$ws.Range("C31").Value = 288
$ws.Range("D31").Value = -172
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.40 : 1"

$ws.Range("B32").Value = 485
$ws.Range("C32").Value = 557
$ws.Range("D32").Value = -72
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.87 : 1"

$ws.Range("B33").Value = 428
$ws.Range("C33").Value = 477
$ws.Range("D33").Value = -49
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.90 : 1"

$ws.Range("B34").Value = 181
$ws.Range("C34").Value = 119
$ws.Range("D34").Value = 62
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.52 : 1"

$ws.Range("B35").Value = 912
$ws.Range("C35").Value = 1041
$ws.Range("D35").Value = -129
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.88 : 1"

$ws.Range("B36").Value = 198
$ws.Range("C36").Value = 497
$ws.Range("D36").Value = -299
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.40 : 1"

$ws.Range("B37").Value = 494
$ws.Range("C37").Value = 333
$ws.Range("D37").Value = 161
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.48 : 1"

$ws.Range("B38").Value = 20
$ws.Range("C38").Value = 177
$ws.Range("D38").Value = -157
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.11 : 1"

$ws.Range("B39").Value = 8
$ws.Range("C39").Value = 64
$ws.Range("D39").Value = -56
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.13 : 1"

$ws.Range("B40").Value = 110
$ws.Range("C40").Value = 107
$ws.Range("D40").Value = 3
$ws.Range("E40").Value = "We borrowerd more than we lent"
$ws.Range("G40").Value = "1.03 : 1"

$ws.Range("B41").Value = 3
$ws.Range("C41").Value = 22
$ws.Range("D41").Value = -19
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.14 : 1"

$ws.Range("B42").Value = 26
$ws.Range("C42").Value = 35
$ws.Range("D42").Value = -9
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.74 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 93
$ws.Range("C44").Value = 76
$ws.Range("D44").Value = 17
$ws.Range("E44").Value = "We borrowerd more than we lent"
$ws.Range("G44").Value = "1.22 : 1"

$ws.Range("B45").Value = 89
$ws.Range("C45").Value = 180
$ws.Range("D45").Value = -91
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.49 : 1"

$ws.Range("B46").Value = 453
$ws.Range("C46").Value = 676
$ws.Range("D46").Value = -223
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.67 : 1"

$ws.Range("B47").Value = 1006
$ws.Range("C47").Value = 697
$ws.Range("D47").Value = 309
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.44 : 1"

$ws.Range("B48").Value = 332
$ws.Range("C48").Value = 543
$ws.Range("D48").Value = -211
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.61 : 1"

$ws.Range("B49").Value = 415
$ws.Range("C49").Value = 266
$ws.Range("D49").Value = 149
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "1.56 : 1"

$ws.Range("B50").Value = 1015
$ws.Range("C50").Value = 491
$ws.Range("D50").Value = 524
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "2.07 : 1"

$ws.Range("B51").Value = 214
$ws.Range("C51").Value = 158
$ws.Range("D51").Value = 56
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.35 : 1"

$ws.Range("B52").Value = 408
$ws.Range("C52").Value = 510
$ws.Range("D52").Value = -102
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.80 : 1"

$ws.Range("B53").Value = 120
$ws.Range("C53").Value = 250
$ws.Range("D53").Value = -130
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.48 : 1"

$ws.Range("B54").Value = 29
$ws.Range("C54").Value = 218
$ws.Range("D54").Value = -189
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.13 : 1"

$ws.Range("B55").Value = 329
$ws.Range("C55").Value = 211
$ws.Range("D55").Value = 118
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "1.56 : 1"
